$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A45").Value = "Report 01 (Review + modify)"
$ws.Range("B45").Value = 10
$ws.Range("C45").Value = "Finish task on time"

$ws.Range("A46").Value = "Report 02 (Review)"
$ws.Range("B46").Value = 10
$ws.Range("C46").Value = "Finish task on time"

$ws.Range("A47").Value = "Report 03 SRS_PMS (Page 19-32, 39-64), User Requirement_PMS(4.Project Eye, 7.Admin)"
$ws.Range("B47").Value = 10
$ws.Range("C47").Value = "Finish task on time"

$ws.Range("A48").Value = "Report 04 ( Page 42-46, 56-74)"
$ws.Range("B48").Value = 10
$ws.Range("C48").Value = "Finish task on time"
